$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add formula in D15: cube root of B15
$ws.Range("D15").Formula = "=(B15)^(1/3)"

# Update B16 formula: B15/B13 instead of SUM(B15)
$ws.Range("B16").Formula = "=B15/B13"
